$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 was a blank placeholder row (just the index 19 and blank formatting).
# Clone the formatting of the row above (row 20) so the new row matches the
# same alternating style/borders used throughout the table, then overwrite
# the values with the new expense entry.
$ws.Range("B20:G20").Copy($ws.Range("B21:G21"))

$ws.Range("B21").Value = 19
$ws.Range("C21").Value = "支出"
$ws.Range("D21").Value = 300
$ws.Range("E21").Value = (Get-Date -Year 2017 -Month 10 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F21").Value = "生活费"
$ws.Range("G21").Value = "生活费(11/01-11/10)"

# Move the selection to M14 (matches the author's saved cursor position).
[void]$ws.Range("M14").Select()
